# Add two new header/value column pairs (I: "I0" / 9, J: "IF" / 9)
# mirroring the existing header style used by column H ("IP").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the existing header cell H1 onto the
# new header cells I1 and J1, so they share the same bold/centered/
# bordered header style without creating duplicate style entries.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Set the header labels for the new columns.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Set the data values for row 2 in the new columns.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
